$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5250
$ws.Range("J62").Value = 6000
$ws.Range("L62").Value = 6000
$ws.Range("N62").Value = -7248

$ws.Range("H65").Value = 5250
$ws.Range("J65").Value = 6000
$ws.Range("L65").Value = 30000
$ws.Range("N65").Value = -36240

$ws.Range("H86").Value = 2710.05
$ws.Range("I86").Value = 2566.8
$ws.Range("J86").Value = 3139.8
$ws.Range("K86").Value = 2566.8
$ws.Range("L86").Value = 3139.8
$ws.Range("M86").Value = -1443.8
$ws.Range("N86").Value = -5385.8

$ws.Range("H89").Value = 2710.05
$ws.Range("I89").Value = 2566.8
$ws.Range("J89").Value = 3139.8
$ws.Range("K89").Value = 12834
$ws.Range("L89").Value = 15699
$ws.Range("M89").Value = -7218
$ws.Range("N89").Value = -26931

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7683.33
$ws.Range("I32").Value = 5567.5444
$ws.Range("J32").Value = 15642.714
$ws.Range("K32").Value = 5567.5444
$ws.Range("L32").Value = 15642.714
$ws.Range("M32").Value = -5280.5444
$ws.Range("N32").Value = -16216.714

$ws.Range("H97").Value = 2994.0667
$ws.Range("I97").Value = 4170
$ws.Range("J97").Value = 642.2
$ws.Range("K97").Value = 4170
$ws.Range("L97").Value = 642.2
$ws.Range("M97").Value = -3674
$ws.Range("N97").Value = -1634.2

$ws.Range("H105").Value = 35500
$ws.Range("J105").Value = 35500
$ws.Range("L105").Value = 35500
$ws.Range("N105").Value = -42488

$ws.Range("H135").Value = 2157945.8
$ws.Range("J135").Value = 2157945.8
$ws.Range("L135").Value = 2157945.8
$ws.Range("N135").Value = -2168085.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 790
$ws.Range("J64").Value = 790
$ws.Range("L64").Value = 790
$ws.Range("N64").Value = -1240

$ws.Range("H67").Value = 790
$ws.Range("J67").Value = 790
$ws.Range("L67").Value = 790
$ws.Range("N67").Value = -2350

$ws.Range("H86").Value = 7695108
$ws.Range("I86").Value = 12502613
$ws.Range("K86").Value = 12502613
$ws.Range("M86").Value = -12501490

$ws.Range("H89").Value = 7695108
$ws.Range("I89").Value = 12502613
$ws.Range("K89").Value = 62513065
$ws.Range("M89").Value = -62507449

$ws.Range("H97").Value = 4539.5557
$ws.Range("I97").Value = 1607
$ws.Range("J97").Value = 28000
$ws.Range("K97").Value = 1607
$ws.Range("L97").Value = 28000
$ws.Range("M97").Value = -616
$ws.Range("N97").Value = -29982

$ws.Range("H134").Value = 251291.88
$ws.Range("I134").Value = 370973.97
$ws.Range("J134").Value = 2721.3845
$ws.Range("K134").Value = 1112921.91
$ws.Range("L134").Value = 8164.1535
$ws.Range("M134").Value = -1110386.91
$ws.Range("N134").Value = -13234.1535

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1655.2616
$ws.Range("I31").Value = 1102.7333
$ws.Range("J31").Value = 2898.45
$ws.Range("K31").Value = 1102.7333
$ws.Range("L31").Value = 2898.45
$ws.Range("M31").Value = -807.7333000000001
$ws.Range("N31").Value = -3488.45

$ws.Range("H34").Value = 1655.2616
$ws.Range("I34").Value = 1102.7333
$ws.Range("J34").Value = 2898.45
$ws.Range("K34").Value = 1102.7333
$ws.Range("L34").Value = 2898.45
$ws.Range("M34").Value = -900.7333000000001
$ws.Range("N34").Value = -3302.45

$ws.Range("H39").Value = 10000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 10000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 10000
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -10782

$ws.Range("H49").Value = 10000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 10000
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 10000
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -10364

$ws.Range("H97").Value = 29660
$ws.Range("J97").Value = 29660
$ws.Range("L97").Value = 29660
$ws.Range("N97").Value = -31642

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 856.7222
$ws.Range("I5").Value = 509.75
$ws.Range("J5").Value = 1290.4375
$ws.Range("K5").Value = 1529.25
$ws.Range("L5").Value = 3871.3125
$ws.Range("M5").Value = -1417.25
$ws.Range("N5").Value = -4095.3125

$ws.Range("H13").Value = 244.42857
$ws.Range("I13").Value = 82.2
$ws.Range("J13").Value = 650
$ws.Range("K13").Value = 246.6
$ws.Range("L13").Value = 1950
$ws.Range("M13").Value = -78.60000000000002
$ws.Range("N13").Value = -2286

$ws.Range("H16").Value = 1171.7142
$ws.Range("I16").Value = 860
$ws.Range("J16").Value = 1951
$ws.Range("K16").Value = 2580
$ws.Range("L16").Value = 5853
$ws.Range("M16").Value = -2407
$ws.Range("N16").Value = -6199

$ws.Range("H20").Value = 1396
$ws.Range("J20").Value = 1450
$ws.Range("L20").Value = 4350
$ws.Range("N20").Value = -4804

$ws.Range("H103").Value = 1919.2667
$ws.Range("I103").Value = 393
$ws.Range("K103").Value = 1179
$ws.Range("M103").Value = -300

$ws.Range("H122").Value = 1134.8235
$ws.Range("I122").Value = 565.6667
$ws.Range("J122").Value = 1445.2727
$ws.Range("K122").Value = 5091.0003
$ws.Range("L122").Value = 13007.4543
$ws.Range("M122").Value = -2641.0003
$ws.Range("N122").Value = -17907.4543

$ws.Range("H135").Value = 856.7222
$ws.Range("I135").Value = 509.75
$ws.Range("J135").Value = 1290.4375
$ws.Range("K135").Value = 4587.75
$ws.Range("L135").Value = 11613.9375
$ws.Range("M135").Value = -2052.75
$ws.Range("N135").Value = -16683.9375

$ws.Range("H136").Value = 1038.35
$ws.Range("I136").Value = 906.0909
$ws.Range("J136").Value = 1200
$ws.Range("K136").Value = 2718.2727
$ws.Range("L136").Value = 3600
$ws.Range("M136").Value = 2381.7273
$ws.Range("N136").Value = -13800

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H113").Value = 1725.7
$ws.Range("I113").Value = 1420
$ws.Range("J113").Value = 1929.5
$ws.Range("K113").Value = 1420
$ws.Range("L113").Value = 1929.5
$ws.Range("M113").Value = 750
$ws.Range("N113").Value = -6269.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 500
$ws.Range("I16").Value = 500
$ws.Range("K16").Value = 500
$ws.Range("M16").Value = -330

$ws.Range("H46").Value = 800.3333
$ws.Range("I46").Value = 650.5
$ws.Range("J46").Value = 1100
$ws.Range("K46").Value = 650.5
$ws.Range("L46").Value = 1100
$ws.Range("M46").Value = -462.5
$ws.Range("N46").Value = -1476

$ws.Range("H55").Value = 300.74194
$ws.Range("I55").Value = 387.92856
$ws.Range("J55").Value = 228.94118
$ws.Range("K55").Value = 387.92856
$ws.Range("L55").Value = 228.94118
$ws.Range("M55").Value = -214.92856
$ws.Range("N55").Value = -574.94118

$ws.Range("H59").Value = 5000
$ws.Range("J59").Value = 5000
$ws.Range("L59").Value = 5000
$ws.Range("N59").Value = -6308

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 29625.5
$ws.Range("J80").Value = 29625.5
$ws.Range("L80").Value = 29625.5
$ws.Range("N80").Value = -31621.5

$ws.Range("H83").Value = 29625.5
$ws.Range("J83").Value = 29625.5
$ws.Range("L83").Value = 88876.5
$ws.Range("N83").Value = -98860.5

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws.Range("H132").Value = 1819.8572
$ws.Range("I132").Value = 827.2105
$ws.Range("J132").Value = 2998.625
$ws.Range("K132").Value = 2481.6315
$ws.Range("L132").Value = 8995.875
$ws.Range("M132").Value = 48.36850000000004
$ws.Range("N132").Value = -14055.875

$ws.Range("H136").Value = 5614.609
$ws.Range("I136").Value = 579.7692
$ws.Range("J136").Value = 12159.9
$ws.Range("K136").Value = 1739.3076
$ws.Range("L136").Value = 36479.7
$ws.Range("M136").Value = 810.6924000000001
$ws.Range("N136").Value = -41579.7
